$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text could be misread as a number by Excel when assigned via .Value
# are first forced to Text format, written, then restored to the Normal style so
# no stray number-format is left behind on the cell.
$textCells = @(
    "D2",
    "D3",
    "D5",
    "D6",
    "D9",
    "D10",
    "D11",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D21",
    "D22",
    "D24",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D51",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '45.525.80'
$ws.Range('E2').Value = '  +7.11%  '
$ws.Range('D3').Value = '2.382.31'
$ws.Range('E3').Value = '  +4.12%  '
$ws.Range('E4').Value = '  +0.79%  '
$ws.Range('D5').Value = '114.58'
$ws.Range('E5').Value = '  +10.81%  '
$ws.Range('D6').Value = '317.37'
$ws.Range('E6').Value = '  +2.13%  '
$ws.Range('E7').Value = '  +2.07%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Value = '0.629'
$ws.Range('E9').Value = '  +4.61%  '
$ws.Range('D10').Value = '42.77'
$ws.Range('E10').Value = '  +10.73%  '
$ws.Range('D11').Value = '0.0938'
$ws.Range('E11').Value = '  +4.15%  '
$ws.Range('E12').Value = '  +6.04%  '
$ws.Range('E13').Value = '  +4.95%  '
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').Value = '15.87'
$ws.Range('E15').Value = '  +4.33%  '
$ws.Range('D16').Value = '2.742.88'
$ws.Range('E16').Value = '  +4.23%  '
$ws.Range('D17').Value = '2.382.66'
$ws.Range('E17').Value = '  +4.46%  '
$ws.Range('D18').Value = '45.403.04'
$ws.Range('E18').Value = '  +6.96%  '
$ws.Range('D19').Value = '7.60'
$ws.Range('E19').Value = '  +4.29%  '
$ws.Range('E20').Value = '  +3.84%  '
$ws.Range('D21').Value = '13.33'
$ws.Range('E21').Value = '  +0.28%  '
$ws.Range('D22').Value = '74.80'
$ws.Range('E22').Value = '  +2.25%  '
$ws.Range('E23').Value = '  +4.23%  '
$ws.Range('D24').Value = '268.81'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('E25').Value = '  +9.24%  '
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('E27').Value = '  +9.48%  '
$ws.Range('E28').Value = '  +5.51%  '
$ws.Range('D29').Value = '2.34'
$ws.Range('E29').Value = '  +2.47%  '
$ws.Range('D30').Value = '22.91'
$ws.Range('E30').Value = '  +2.63%  '
$ws.Range('D31').Value = '38.86'
$ws.Range('E31').Value = '  +8.44%  '
$ws.Range('D32').Value = '0.0961'
$ws.Range('E32').Value = '  +13.65%  '
$ws.Range('D33').Value = '170.66'
$ws.Range('E33').Value = '  +3.87%  '
$ws.Range('E34').Value = '  +16.94%  '
$ws.Range('E35').Value = '  +1.84%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '4.98'
$ws.Range('E36').Value = '  +10.84%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '0.120'
$ws.Range('E37').Value = '  +7.75%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '4.10'
$ws.Range('E38').Value = '  +13.60%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '3.07'
$ws.Range('E39').Value = '  +11.69%  '
$ws.Range('D40').Value = '0.0366'
$ws.Range('E40').Value = '  +6.35%  '
$ws.Range('E41').Value = '  +11.08%  '
$ws.Range('D42').Value = '104.34'
$ws.Range('E42').Value = '  -7.28%  '
$ws.Range('E43').Value = '  +6.55%  '
$ws.Range('D44').Value = '71.39'
$ws.Range('E44').Value = '  +2.57%  '
$ws.Range('D45').Value = '13.33'
$ws.Range('E45').Value = '  +10.32%  '
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').Value = '5.77'
$ws.Range('E47').Value = '  +12.04%  '
$ws.Range('D48').Value = '116.21'
$ws.Range('E48').Value = '  +5.55%  '
$ws.Range('D49').Value = '1.65'
$ws.Range('E49').Value = '  +18.10%  '
$ws.Range('E50').Value = '  +8.30%  '
$ws.Range('D51').Value = '79.46'
$ws.Range('E51').Value = '  +3.26%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
